# Actualización desde MV -datos-
# Adds three new daily rows (08-09-2021, 09-09-2021, 10-09-2021) to the
# bottom of the data table on Sheet1, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("08-09-2021", 5120, 629, 695, 1115, 1281, 1400),
    @("09-09-2021", 5869, 913, 856, 1603, 1207, 1289),
    @("10-09-2021", 4392, 602, 697, 1249, 906, 938)
)

$startRow = 174
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A holds a date-formatted label (e.g. "08-09-2021"). Assigning a
    # plain string would let Excel auto-convert it to a date serial; force
    # text interpretation, then restore the default "Normal" style so the
    # new cell matches the formatting of the existing rows (no explicit
    # style index).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
